$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2022-03-16"

# Update the label for the March row in column A
$ws.Range("A4").Value = "March (through 03-16)"

# Update March row (row 4) values
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 31
$ws.Range("E4").Value = 33
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = 71

# Update Total row (row 5) values
$ws.Range("B5").Value = 52
$ws.Range("C5").Value = 107
$ws.Range("D5").Value = 162
$ws.Range("E5").Value = 170
$ws.Range("H5").Value = 387
$ws.Range("I5").Value = 371
